$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coinranking price/volume table refresh (scheduled GitHub Actions run).
# D (Price) and E (Volume 1h) columns store plain text in the source sheet,
# so each target cell is forced to Text format before the value is written --
# otherwise numeric-looking strings such as "19.63" or "2.000" would be
# re-interpreted as numbers and lose significant trailing zeros / formatting.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.908.44"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -2.17%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.833.78"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.68%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.70%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4614"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -1.32%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3666"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.80%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07177"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -2.78%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8801"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.13%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07841"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.45%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "19.63"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.99%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.852.68"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.44%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.335"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.70%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.393"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -3.38%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "88.51"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -4.62%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008765"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.92%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.02%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "26.929.29"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -2.18%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -2.88%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.017"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -2.80%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.43"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.28%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.973"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +5.34%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "150.80"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.70%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "18.24"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -1.65%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.000"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -4.61%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "113.51"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -3.12%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.962"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -4.01%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08845"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.80%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.133"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +3.83%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7692"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +1.88%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.469"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.61%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.135"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -2.61%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.658"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.03%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.091"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.73%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01935"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.95%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.929"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -2.08%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05146"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -2.75%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.962"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -3.15%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4980"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -4.67%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1599"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -3.01%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.346"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.33%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4695"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -3.88%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.23"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.59%  "
$ws.Range("B46").Value = "PaxDollar"
$ws.Range("C46").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.004"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.06%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "103.10"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.71%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.617"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -2.85%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06101"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -2.59%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "64.94"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.60%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "36.48"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -2.02%  "
